# Icons.xlsx update — add "Xuất lắp ráp" (assembly export) menu/icon entry
# and backfill the MODULE / GUI columns that were missing on the two
# preceding rows of the "Phiếu xuất kho" (StockInOut) block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29 : fill in the MODULE (A) / GUI (B) columns that were blank ---
$ws.Range("A29").Value = "VnsErp2025"
$ws.Range("B29").Value = "FormMain"

# --- Row 30 : fill in MODULE (A) / GUI (B) / PHÂN LOẠI (C) that were blank ---
$ws.Range("A30").Value = "VnsErp2025"
$ws.Range("B30").Value = "FormMain"
$ws.Range("C30").Value = "StockInOutRibbonPageGroup"

# --- Row 31 (new row) : new "Xuất lắp ráp" assembly-export icon entry ---
$ws.Range("A31").Value = "VnsErp2025"
$ws.Range("B31").Value = "FormMain"
$ws.Range("C31").Value = "XuatKhoRibbonPageGroup"
$ws.Range("D31").Value = "XuatLapRapBarButtonItem"
$ws.Range("E31").Value = "Xuất lắp ráp"
$ws.Range("F31").Value = "dashboard.svg"

# Move / resize the selection onto the newly added row, matching the
# author's saved cursor position.
$ws.Range("A31:B31").Select()
